$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (main schedule sheet)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Range("A2").Value = "Última actualización: 10:49:38"
$ws1.Range("A3").Value = "Total filas: 116"

# Swap the data in rows 62 and 63 (Hora_Llegada / Parada stay the same,
# Hora_Scrap / Linea / Minutos trade places between the two rows)
$ws1.Range("A62").Value = "08:38:24"
$ws1.Range("C62").Value = "27_EL RETIRO"
$ws1.Range("D62").Value = 39

$ws1.Range("A63").Value = "07:49:32"
$ws1.Range("C63").Value = "14_ABASTO"
$ws1.Range("D63").Value = 88

# Insert a new row at position 109 (pushes old 109-118 down to 110-119)
$ws1.Rows.Item(109).Insert()
$ws1.Range("A109").Value = "10:49:38"
$ws1.Range("B109").Value = "11:47"
$ws1.Range("C109").Value = "23_HERNANDEZ"
$ws1.Range("D109").Value = 58
$ws1.Range("E109").Value = "LP1912"

# Append two brand-new rows at the end (120 and 121)
$ws1.Range("A120").Value = "10:49:38"
$ws1.Range("B120").Value = "12:36"
$ws1.Range("C120").Value = "27_EL RETIRO"
$ws1.Range("D120").Value = 107
$ws1.Range("E120").Value = "LP1912"

$ws1.Range("A121").Value = "10:49:38"
$ws1.Range("B121").Value = "12:48"
$ws1.Range("C121").Value = "16_SANTA ANA"
$ws1.Range("D121").Value = 119
$ws1.Range("E121").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 10:49:38"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 10:49:38"
